# Update the evaluate script
#
# The "OpenMPI" row in each results table is replaced by three separate
# rows ("OpenMPI - 2 thread", "OpenMPI - 4 thread", "OpenMPI - 8 thread"),
# applied identically to all three worksheets (a1a, ijcnn1, generated).
# The previous OpenMPI row (row 5) is removed - which shifts the CUDA row
# (previously row 6) up into row 5 - and three new label-only rows are
# appended at rows 7-9.

$wb = $excel.ActiveWorkbook

$threadLabels = @("OpenMPI - 2 thread", "OpenMPI - 4 thread", "OpenMPI - 8 thread")

for ($i = 1; $i -le 3; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Remove the old "OpenMPI" row; this shifts all following rows up by one
    # (the old CUDA row 6 becomes row 5).
    $ws.Rows.Item(5).Delete()

    # Clear out the target area for the new rows, then write just the labels
    # in column A (rows 7-9), leaving the rest of the row blank, as in the diff.
    $ws.Range("A7:H9").Value = $null
    $ws.Range("A7").Value = $threadLabels[0]
    $ws.Range("A8").Value = $threadLabels[1]
    $ws.Range("A9").Value = $threadLabels[2]
}

# Update the selections on each sheet; doing the "generated" sheet (sheet3)
# last keeps it the active tab, matching the original workbook.
[void]$wb.Worksheets.Item(1).Range("A7:A9").Select()
[void]$wb.Worksheets.Item(2).Range("A7:A9").Select()
[void]$wb.Worksheets.Item(3).Range("A7:A9").Select()
